# Example.xlsx -> rename/add sheets for the "Linear" regression example
# and fix up the theme's Light-1 color (c0c0c0 -> white), per the
# "Linear eq introduced and tested: testLinear.py" commit.

$wb = $excel.ActiveWorkbook

# 1) Rename the first sheet "Example" -> "Dose_Resp" (keeps sheetId/rId).
$wsDoseResp = $wb.Worksheets.Item(1)
$wsDoseResp.Name = "Dose_Resp"

# 2) Append a brand-new worksheet named "Linear" as the last (3rd) tab.
#    Adding relative to the last existing sheet places it at the end,
#    matching sheetId="3" / r:id="rId3" / sheet3.xml ordering, and makes
#    it the active sheet (workbookView activeTab=2, sheetView tabSelected=1).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLinear = $wb.Worksheets.Add($null, $lastSheet)
$wsLinear.Name = "Linear"

# 3) Populate the Linear sheet with the dose/response-style paired data.
$linearData = @(
    @(1, 3.17),
    @(1, 13.25),
    @(2, 19.8),
    @(2, 14.18),
    @(3, 11.43),
    @(3, 25.85),
    @(4, 13.81),
    @(4, 25.49),
    @(5, 26.94),
    @(5, 38.86)
)

for ($i = 0; $i -lt $linearData.Count; $i++) {
    $row = $i + 1
    $wsLinear.Cells.Item($row, 1).Value = $linearData[$i][0]
    $wsLinear.Cells.Item($row, 2).Value = $linearData[$i][1]
}

# Leave the selection on the last data cell, as in the authored sheet.
$wsLinear.Range("A10").Select()

# 4) The workbook theme's Light-1 color changed from C0C0C0 to pure white.
$themeColors = $wb.Theme.ThemeColorScheme
$themeColors.Colors(2).RGB = 0xFFFFFF
